$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-07-03 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-07-04 Thursday", 2) | Out-Null
$d.Content.Find.Execute("436×2=872", $true, $false, $false, $false, $false, $true, 1, $false, "320×6=1920", 2) | Out-Null
$d.Content.Find.Execute("937×3=2811", $true, $false, $false, $false, $false, $true, 1, $false, "165×3=495", 2) | Out-Null
$d.Content.Find.Execute("303×4=1212", $true, $false, $false, $false, $false, $true, 1, $false, "316×2=632", 2) | Out-Null
$d.Content.Find.Execute("736×9=6624", $true, $false, $false, $false, $false, $true, 1, $false, "437×7=3059", 2) | Out-Null
$d.Content.Find.Execute("638×2=1276", $true, $false, $false, $false, $false, $true, 1, $false, "546×2=1092", 2) | Out-Null
$d.Content.Find.Execute("603×6=3618", $true, $false, $false, $false, $false, $true, 1, $false, "458×9=4122", 2) | Out-Null
$d.Content.Find.Execute("741×2=1482", $true, $false, $false, $false, $false, $true, 1, $false, "186×6=1116", 2) | Out-Null
$d.Content.Find.Execute("486×2=972", $true, $false, $false, $false, $false, $true, 1, $false, "606×5=3030", 2) | Out-Null
$d.Content.Find.Execute("587×4=2348", $true, $false, $false, $false, $false, $true, 1, $false, "688×8=5504", 2) | Out-Null
$d.Content.Find.Execute("221×5=1105", $true, $false, $false, $false, $false, $true, 1, $false, "168×7=1176", 2) | Out-Null
$d.Content.Find.Execute("875×8=7000", $true, $false, $false, $false, $false, $true, 1, $false, "630×8=5040", 2) | Out-Null
$d.Content.Find.Execute("805×5=4025", $true, $false, $false, $false, $false, $true, 1, $false, "685×5=3425", 2) | Out-Null
$d.Content.Find.Execute("925×3=2775", $true, $false, $false, $false, $false, $true, 1, $false, "262×7=1834", 2) | Out-Null
$d.Content.Find.Execute("719×5=3595", $true, $false, $false, $false, $false, $true, 1, $false, "117×7=819", 2) | Out-Null
$d.Content.Find.Execute("802×4=3208", $true, $false, $false, $false, $false, $true, 1, $false, "384×2=768", 2) | Out-Null
$d.Content.Find.Execute("843×7=5901", $true, $false, $false, $false, $false, $true, 1, $false, "584×9=5256", 2) | Out-Null
$d.Content.Find.Execute("801×2=1602", $true, $false, $false, $false, $false, $true, 1, $false, "230×3=690", 2) | Out-Null
$d.Content.Find.Execute("109×7=763", $true, $false, $false, $false, $false, $true, 1, $false, "567×7=3969", 2) | Out-Null
$d.Content.Find.Execute("798×3=2394", $true, $false, $false, $false, $false, $true, 1, $false, "984×9=8856", 2) | Out-Null
$d.Content.Find.Execute("708×7=4956", $true, $false, $false, $false, $false, $true, 1, $false, "679×9=6111", 2) | Out-Null
$d.Content.Find.Execute("743×2=1486", $true, $false, $false, $false, $false, $true, 1, $false, "243×4=972", 2) | Out-Null
$d.Content.Find.Execute("343×2=686", $true, $false, $false, $false, $false, $true, 1, $false, "667×6=4002", 2) | Out-Null
$d.Content.Find.Execute("683×6=4098", $true, $false, $false, $false, $false, $true, 1, $false, "342×2=684", 2) | Out-Null
$d.Content.Find.Execute("547×8=4376", $true, $false, $false, $false, $false, $true, 1, $false, "934×5=4670", 2) | Out-Null
$d.Content.Find.Execute("384×9=3456", $true, $false, $false, $false, $false, $true, 1, $false, "394×5=1970", 2) | Out-Null

Write-Host "Replacements complete"
